# Hortaliza, Feria Lagunitas de Puerto Montt - Pimiento
# Weekly update: insert two new price rows (week of 2022-06-02, serial 44714)
# above the former row 665, pushing the existing rows 665:711 down to 667:713.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 665 (shifts old 665:711 -> 667:713).
$ws.Rows("665:666").Insert()

# New row 665: Zafiro rojo, Primera
$ws.Range("A665").Value = 4
$ws.Range("B665").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C665").Value = "Los Lagos"
$ws.Range("D665").Value = 44714
$ws.Range("E665").Value = 10
$ws.Range("F665").Value = 100112002
$ws.Range("G665").Value = "Pimiento"
$ws.Range("H665").Value = "Zafiro rojo"
$ws.Range("I665").Value = "Primera"
$ws.Range("J665").Value = 90
$ws.Range("K665").Value = 48000
$ws.Range("L665").Value = 48000
$ws.Range("M665").Value = 48000
$ws.Range("N665").Value = "$/caja 15 kilos"
$ws.Range("O665").Value = "Región de Arica y Parinacota"
$ws.Range("P665").Value = 3200
$ws.Range("Q665").Value = 15
$ws.Range("R665").Value = "Hortaliza"

# New row 666: Zafiro verde, Extra
$ws.Range("A666").Value = 4
$ws.Range("B666").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C666").Value = "Los Lagos"
$ws.Range("D666").Value = 44714
$ws.Range("E666").Value = 10
$ws.Range("F666").Value = 100112002
$ws.Range("G666").Value = "Pimiento"
$ws.Range("H666").Value = "Zafiro verde"
$ws.Range("I666").Value = "Extra"
$ws.Range("J666").Value = 90
$ws.Range("K666").Value = 30000
$ws.Range("L666").Value = 30000
$ws.Range("M666").Value = 30000
$ws.Range("N666").Value = "$/caja 15 kilos"
$ws.Range("O666").Value = "Región de Arica y Parinacota"
$ws.Range("P666").Value = 2000
$ws.Range("Q666").Value = 15
$ws.Range("R666").Value = "Hortaliza"
